# Covid-19 data update: append one new daily row (27 May 2020) to the
# "Tabela1" table on the single worksheet, growing it from A1:J77 to
# A1:J78 (dimension, autofilter and table range all expand accordingly).
#
# Data mirrors the previous day's row: Date 43978 (=27/5/2020), Tested(all)
# 77210, Tested(daily) 631, Positive(all) 1473, Positive(daily) 2,
# Hospitalized 7, ICU 2, Discharged 1, Deaths(all) 108, Deaths(daily) 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last data row (77) down into the new row 78 first, so the
# new row inherits row 77's exact cell styles (s="20"/"21"/"22" - the date
# format, the thousands-format, and the general-number format used by the
# most recent stretch of the table) instead of picking up a generic default
# style the way a plain Value assignment into a previously-empty row would.
$ws.Rows("77").Copy()
$ws.Rows("78").Insert(-4121, 0)   # xlShiftDown, no special CopyOrigin

# Fill in the new day's figures.
$ws.Range("A78").Value = 43978
$ws.Range("B78").Value = 77210
$ws.Range("C78").Value = 631
$ws.Range("D78").Value = 1473
$ws.Range("E78").Value = 2
$ws.Range("F78").Value = 7
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 1
$ws.Range("I78").Value = 108
$ws.Range("J78").Value = 0

# Grow the table (ListObject) + its AutoFilter range to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J78"))

# Match the workbook's last-saved selection (the new last row is selected).
$ws.Range("A78:J78").Select() | Out-Null
